# Update "想去人数" (interest count, column F) values on the "展览" (sheet1)
# and "全部类型" (sheet4) worksheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" is the 1st worksheet (rId1 -> sheet1.xml)
$wsExhibition = $wb.Worksheets.Item(1)
$exhibitionUpdates = @{
    2  = 8096
    3  = 119
    4  = 94
    5  = 30773
    7  = 604
    8  = 701
    15 = 390
    17 = 557
    21 = 1114
    23 = 703
    24 = 2338
    25 = 829
    26 = 67
    27 = 1086
    29 = 634
    30 = 1075
}
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" is the 4th worksheet (rId4 -> sheet4.xml)
$wsAll = $wb.Worksheets.Item(4)
$allUpdates = @{
    3  = 8096
    4  = 119
    5  = 94
    7  = 30773
    9  = 604
    10 = 701
    21 = 390
    27 = 557
    31 = 1114
    33 = 703
    34 = 2338
    35 = 829
    36 = 67
    37 = 1086
    40 = 634
    41 = 1075
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
